# Auto-update: PANELES GLOBALES TRIMESTRALES
# Appends 20 additional quarterly-date rows (2026-01-01 .. 2030-10-01) to the
# existing date series in column A of all three worksheets ("nivel",
# "trimestrales", "i.a."), extending each sheet's data from row 225 to row 245.

$wb = $excel.ActiveWorkbook

# Quarter-start date serials (Excel 1900 date system) continuing the
# existing A-column series, which already ends at row 225 with 45931
# (2025-10-01).
$newDates = @(46023, 46113, 46204, 46296, 46388, 46478, 46569, 46661, 46753, 46844, 46935, 47027, 47119, 47209, 47300, 47392, 47484, 47574, 47665, 47757)

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    $startRow = 226
    for ($i = 0; $i -lt $newDates.Length; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, 1).Value = $newDates[$i]
    }
}
